$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump version number ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.1"

# --- New Table sheet: fill in example values on the data row (row 2) ---
$ws = $wb.Worksheets.Item("New Table")

$ws.Range("B2").Value = "NMR spectroscopy assay"
$ws.Range("C2").Value = "OBI"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/OBI_0000623"

$ws.Range("F2").Value = "Varian UNITY INOVA spectrometer"
$ws.Range("G2").Value = "OBI"
$ws.Range("H2").Value = "http://purl.obolibrary.org/obo/OBI_0000558"

$ws.Range("I2").Value = "5 mm CPTCI 1H-13C/15N/D Z-GRD"

# Number of scans is a numeric-looking example value that must stay text,
# matching the template's convention of storing everything as strings.
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "128"
$ws.Range("L2").ClearFormats()

$ws.Range("O2").Value = "1D 1H with presaturation (presat)"

# Magnetic field strength is likewise numeric-looking but must stay text.
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "11.7"
$ws.Range("R2").ClearFormats()

$ws.Range("S2").Value = "tesla"
$ws.Range("T2").Value = "UO"
$ws.Range("U2").Value = "http://purl.obolibrary.org/obo/UO_0000228"
